$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# ---- Header text updates (rich-text runs collapsed to plain text; same visual formatting) ----
$ws.Range("A8").Value = "Volume 31   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/9/2024  Through  12/15/2024"

# ---- Cells whose TYPE flips between number and text (need format+value paste trick to avoid
#      the engine auto-coercing digit-looking text back into a number) ----

# D14: number 1 -> text "0" (copy style+value from C14, which already is text "0" style s13)
$ws.Range("C14").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").PasteSpecial(-4163)

# E14: number -100 -> text "***.*" (copy style+value from E15, already text "***.*" style s13)
$ws.Range("E15").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").PasteSpecial(-4163)

# C22: text "0" -> number 1 (copy number style from D22, then set numeric value)
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1

# C23: number 1 -> text "0" (copy style+value from D23, already text "0" style s13)
$ws.Range("D23").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").PasteSpecial(-4163)

# C31: text "0" -> number 1 (copy number style from F31, then set numeric value)
$ws.Range("F31").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("C31").Value = 1

# D31: text "0" -> number 1 (copy number style from G31, then set numeric value)
$ws.Range("G31").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D31").Value = 1

# E31: text "***.*" -> number 0 (copy style+value from H31, already number 0 style s15)
$ws.Range("H31").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = 0

# ---- Plain numeric updates (style/type unchanged) ----

$ws.Range("L14").Value = 20
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("M15").Value = -27.272727272727
$ws.Range("N15").Value = -81.818181818181
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 15
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 218
$ws.Range("J16").Value = 205
$ws.Range("K16").Value = 6.341463414634
$ws.Range("L16").Value = -17.735849056603
$ws.Range("M16").Value = -19.85294117647
$ws.Range("N16").Value = -83.021806853582
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = -14.285714285714
$ws.Range("I17").Value = 335
$ws.Range("J17").Value = 331
$ws.Range("K17").Value = 1.208459214501
$ws.Range("L17").Value = 15.120274914089
$ws.Range("M17").Value = 56.542056074766
$ws.Range("N17").Value = -59.245742092457
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -55
$ws.Range("I18").Value = 117
$ws.Range("J18").Value = 146
$ws.Range("K18").Value = -19.86301369863
$ws.Range("L18").Value = -18.181818181818
$ws.Range("M18").Value = -34.269662921348
$ws.Range("N18").Value = -93.806246691371
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -36.842105263157
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 58
$ws.Range("H19").Value = -22.413793103448
$ws.Range("I19").Value = 611
$ws.Range("J19").Value = 545
$ws.Range("K19").Value = 12.110091743119
$ws.Range("L19").Value = 7.005253940455
$ws.Range("M19").Value = 64.690026954177
$ws.Range("N19").Value = -46.869565217391
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -18.181818181818
$ws.Range("J20").Value = 215
$ws.Range("K20").Value = -26.046511627907
$ws.Range("L20").Value = -44.791666666666
$ws.Range("M20").Value = 80.681818181818
$ws.Range("N20").Value = -91.016949152542
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = -47.5
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 129
$ws.Range("H21").Value = -19.37984496124
$ws.Range("I21").Value = 1462
$ws.Range("J21").Value = 1459
$ws.Range("K21").Value = 0.205620287868
$ws.Range("L21").Value = -7.818411097099
$ws.Range("M21").Value = 27.241079199303
$ws.Range("N21").Value = -79.271232099815
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("I22").Value = 28
$ws.Range("J22").Value = 41
$ws.Range("K22").Value = -31.70731707317
$ws.Range("L22").Value = -6.666666666666
$ws.Range("M22").Value = 7.692307692307
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = 15.151515151515
$ws.Range("F24").Value = 123
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 1424
$ws.Range("J24").Value = 1363
$ws.Range("K24").Value = 4.475421863536
$ws.Range("L24").Value = -0.070175438596
$ws.Range("M24").Value = 131.921824104235
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 59
$ws.Range("G25").Value = 53
$ws.Range("H25").Value = 11.320754716981
$ws.Range("I25").Value = 676
$ws.Range("J25").Value = 619
$ws.Range("K25").Value = 9.208400646203
$ws.Range("L25").Value = -4.384724186704
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 47
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 62.068965517241
$ws.Range("I26").Value = 554
$ws.Range("J26").Value = 513
$ws.Range("K26").Value = 7.992202729044
$ws.Range("L26").Value = 25.623582766439
$ws.Range("M26").Value = 5.725190839694
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 59
$ws.Range("J28").Value = 54
$ws.Range("K28").Value = 9.259259259259
$ws.Range("L28").Value = 13.461538461538
$ws.Range("D29").Value = 2
$ws.Range("G29").Value = 3
$ws.Range("J29").Value = 17
$ws.Range("K29").Value = 17.647058823529
$ws.Range("L29").Value = -4.761904761904
$ws.Range("N29").Value = -88.304093567251
$ws.Range("G30").Value = 2
$ws.Range("J30").Value = 14
$ws.Range("K30").Value = 21.428571428571
$ws.Range("L30").Value = -19.047619047619
$ws.Range("N30").Value = -88.590604026845
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 2
$ws.Range("I31").Value = 12
$ws.Range("J31").Value = 7
$ws.Range("K31").Value = 71.428571428571
$ws.Range("L31").Value = 140